$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header: "time_taken", styled like the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$xlPasteFormats = -4122
$ws.Range("F1").PasteSpecial($xlPasteFormats)

# time_taken values for each data row (rows 2-14), plain/unstyled like the other data cells
$times = @(
    "2021-10-05 10:51:09.278119",
    "2021-10-05 10:51:09.278129",
    "2021-10-05 10:51:09.278133",
    "2021-10-05 10:51:09.278135",
    "2021-10-05 10:51:09.278138",
    "2021-10-05 10:51:09.278141",
    "2021-10-05 10:51:09.278144",
    "2021-10-05 10:51:09.278146",
    "2021-10-05 10:51:09.278149",
    "2021-10-05 10:51:09.278152",
    "2021-10-05 10:51:09.278154",
    "2021-10-05 10:51:09.278157",
    "2021-10-05 10:51:09.278159"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
